$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.170.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.606.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.66%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("E5").Value = '  +0.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '302.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.32%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3754'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.61%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3631'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.60'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.002'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.262'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08046'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.86'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.538'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.77%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.655'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.98%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001262'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.607.14'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.76%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.30'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.37%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06782'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.562'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.37%  '

$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.72%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.173.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.345'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.80%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.896'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.259'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.92%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.83%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.415'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.716'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -14.18%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.784.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9705'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.40%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07710'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02767'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.214'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2535'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.66%  '

$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.08840'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.24%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.69%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.393'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7129'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.65%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.95'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6561'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.289'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.95%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.978'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.68%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07994'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.44%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.163'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.66%  '

